# Update from MV -datos- : revise the last existing quarter row (74) and
# append the new quarter row (75, "01-04-2021") to the "Por sector
# acreedor" quarterly table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing row 74 values ------------------------------------
$ws.Range("B74").Value = 209510
$ws.Range("C74").Value = 2829
$ws.Range("E74").Value = 2655
$ws.Range("F74").Value = 5255
$ws.Range("H74").Value = 5253
$ws.Range("K74").Value = 33354
$ws.Range("L74").Value = 8504
$ws.Range("M74").Value = 24850
$ws.Range("N74").Value = 72606
$ws.Range("O74").Value = 12382
$ws.Range("P74").Value = 60224
$ws.Range("Q74").Value = 94308
$ws.Range("S74").Value = 92249

# --- Append the new row 75 ---------------------------------------------
# Column A holds period labels as text (e.g. "01-01-2021"). Assigning a
# literal "01-04-2021" string to a cell makes Excel auto-convert it to a
# date serial, so we build it as a text formula result in a scratch cell
# and paste-special (values only) into A75 - this keeps the cell a plain
# shared-string (t="s") without forcing a number-format style on it.
$ws.Range("Z1").Formula = "=""01-04-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

$ws.Range("B75").Value = 213283
$ws.Range("C75").Value = 2664
$ws.Range("D75").Value = 73
$ws.Range("E75").Value = 2591
$ws.Range("F75").Value = 4955
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 4954
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 32910
$ws.Range("L75").Value = 8785
$ws.Range("M75").Value = 24125
$ws.Range("N75").Value = 71901
$ws.Range("O75").Value = 12272
$ws.Range("P75").Value = 59629
$ws.Range("Q75").Value = 99687
$ws.Range("R75").Value = 1186
$ws.Range("S75").Value = 98501
$ws.Range("T75").Value = 1165
$ws.Range("U75").Value = 1165
